# Updated cryptos list (GitHub Actions scrape refresh).
# Price (column D) cells that look like a plain number are forced back to
# Text via NumberFormat "@" before the assignment (otherwise COM/Excel
# auto-converts them to a Number and trailing zeros / exact formatting is
# lost), then the style is reset to "Normal" so no stray number-format
# style is left attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.779.72'
$ws.Range("E2").Value = '  -1.89%  '
$ws.Range("D3").Value = '1.548.36'
$ws.Range("E3").Value = '  -1.86%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '204.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.72%  '
$ws.Range("E6").Value = '  -1.73%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.245'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.10%  '
$ws.Range("B9").Value = 'Solana'
$ws.Range("C9").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.37'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.10%  '
$ws.Range("E10").Value = '  -1.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0858'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.94%  '
$ws.Range("D12").Value = '1.768.18'
$ws.Range("E12").Value = '  -1.87%  '
$ws.Range("D13").Value = '1.552.69'
$ws.Range("E13").Value = '  -1.54%  '
$ws.Range("E14").Value = '  -2.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.511'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.90%  '
$ws.Range("D16").Value = '26.783.10'
$ws.Range("E16").Value = '  -1.89%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.89'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '213.84'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.0₃0683'
$ws.Range("E19").Value = '  -0.81%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.78%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.40%  '
$ws.Range("E24").Value = '  -1.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.72'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.58%  '
$ws.Range("E26").Value = '  -2.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.87'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("E29").Value = '  -2.05%  '
$ws.Range("E30").Value = '  -0.85%  '
$ws.Range("E31").Value = '  -3.87%  '
$ws.Range("E32").Value = '  -0.49%  '
$ws.Range("D33").Value = '1.352.01'
$ws.Range("E33").Value = '  -4.22%  '
$ws.Range("E34").Value = '  -1.00%  '
$ws.Range("E35").Value = '  -3.84%  '
$ws.Range("E36").Value = '  -0.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.917'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.43%  '
$ws.Range("E38").Value = '  -2.19%  '
$ws.Range("E39").Value = '  +0.62%  '
$ws.Range("E40").Value = '  -2.52%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.56'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.07%  '
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.990'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.18'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("E45").Value = '  -2.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '62.85'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.29'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.65%  '
$ws.Range("D48").Value = '1.682.37'
$ws.Range("E48").Value = '  -1.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.81'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0506'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.47%  '
$ws.Range("D51").Value = '0.0₇0971'
$ws.Range("E51").Value = '  -2.04%  '
